$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows of data at the bottom of the table
$ws.Range("A82").Value = 3370
$ws.Range("B82").Value = "Smallest Number With All Set Bits"
$ws.Range("C82").Value = "Math"
$ws.Range("D82").Value = "Power of 2 minus 1"

$ws.Range("A83").Value = 1526
$ws.Range("B83").Value = "Minimum Number of Increments on Subarrays to Form a Target Array"
$ws.Range("C83").Value = "Math"

# Update the view to reflect where the user ended up after the edit
$ws.Range("C69").Select()
$excel.ActiveWindow.ScrollRow = 64
$excel.ActiveWindow.ScrollColumn = 1
